$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns C:F ---
$ws.Range("C1").Value = "Keywords"
$ws.Range("D1").Value = "collection window"
$ws.Range("E1").Value = "population"
$ws.Range("F1").Value = "sample type"

# --- Row 2 (Oxford) ---
$ws.Range("C2").Value = "colostrum; hormones; early postpartum"
$ws.Range("D2").Value = "Days 1-5 Postpartum & 24days+ "
$ws.Range("E2").Value = "term infants"
$ws.Range("F2").Value = "colostrum, transitional, mature milk"

# --- Row 3 (Brooklyn) ---
$ws.Range("C3").Value = "supplements"
$ws.Range("D3").Value = "varied"

# --- Row 4 (NeoBANK) ---
$ws.Range("C4").Value = "NICU; donor milk; preterm"
$ws.Range("D4").Value = "NICU stay"
$ws.Range("E4").Value = "preterm NICU infants"

# --- Row 5 (new study: DHM Pooled) ---
$ws.Range("A5").Value = "DHM Pooled"
$ws.Range("B5").Value = "Looking at single donor profiles and paired samples (pre-post)"
$ws.Range("C5").Value = "Milk Banks"
$ws.Range("D5").Value = "varied"
$ws.Range("E5").Value = "donor milk only"
$ws.Range("F5").Value = "donor human milk (pooled)"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 75.66666666666667
$ws.Columns.Item(3).ColumnWidth = 37.5
$ws.Columns.Item(4).ColumnWidth = 26.5

# --- Selection as left by the author ---
$ws.Range("G7").Select() | Out-Null
